$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: plateWellIDs label becomes the new "noSoln_modelCell" condition; ratio normalized to 1
$ws.Range("B2").Value = "noSoln_modelCell_d1_chamber1_channel1_freq"
$ws.Range("K2").Value = 1

# Row 3: label shifts to the condition previously in row 2
$ws.Range("B3").Value = "16HBEmedia_8TW_d1_chamber1_channel1_freq"

# Row 4 gets real data now (style cloned from the rows above, values set after)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "16HBEmedia_8TWclean_d1_chamber1_channel1_freq"

$ws.Range("K4").Value = 1.1200000000000001

# Rows 5 and 6 (previously empty placeholder rows) are removed
$ws.Rows("5:6").Delete()

# Restore the saved selection/view
$ws.Range("F18").Select()
